$wb = $excel.ActiveWorkbook

# Rename the sheet "This is a sheet" to "This is a - sheet".
# Renaming via the Excel object model automatically updates any formulas
# (e.g. in Sheet2) that reference the sheet via a single-quoted name.
$ws1 = $wb.Worksheets.Item("This is a sheet")
$ws1.Name = "This is a - sheet"
$ws1.Activate()
